$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings. Some look like plain decimals (e.g. "313.03")
# which Excel would auto-convert to a Number on assignment; others use a
# "thousands.thousands" display format (e.g. "42.881.25") that is not a valid
# number so it is kept as text automatically. Force the whole column to Text
# first so every assignment below is stored as a literal string, matching the
# original inline-string cells, then restore the default (unstyled) look.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "42.881.25"
$ws.Range("D3").Value = "2.572.64"
$ws.Range("D5").Value = "313.03"
$ws.Range("D6").Value = "99.34"
$ws.Range("D9").Value = "0.536"
$ws.Range("D10").Value = "35.95"
$ws.Range("D11").Value = "0.0812"
$ws.Range("D12").Value = "7.48"
$ws.Range("D13").Value = "2.967.33"
$ws.Range("D15").Value = "15.92"
$ws.Range("D16").Value = "2.609.01"
$ws.Range("D18").Value = "42.944.59"
$ws.Range("D20").Value = "12.59"
$ws.Range("D22").Value = "69.77"
$ws.Range("D23").Value = "250.36"
$ws.Range("D24").Value = "2.95"
$ws.Range("D25").Value = "2.07"
$ws.Range("D29").Value = "39.81"
$ws.Range("D31").Value = "158.82"
$ws.Range("D32").Value = "5.82"
$ws.Range("D33").Value = "3.36"
$ws.Range("D35").Value = "0.0801"
$ws.Range("D37").Value = "18.63"
$ws.Range("D41").Value = "23.30"
$ws.Range("D42").Value = "4.15"
$ws.Range("D46").Value = "2.003.90"
$ws.Range("D47").Value = "9.05"
$ws.Range("D48").Value = "2.818.61"
$ws.Range("D49").Value = "0.197"
$ws.Range("D50").Value = "82.05"
$ws.Range("D51").Value = "74.75"

$priceRange.Style = "Normal"

# Column E volume strings already carry surrounding spaces + a "%" sign, so
# they are never misread as numbers and can be assigned directly.
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("E15").Value = "  +4.60%  "
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("E20").Value = "  -4.06%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("E29").Value = "  -2.61%  "
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("E32").Value = "  -2.20%  "
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("E35").Value = "  +2.19%  "
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("E38").Value = "  +11.21%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("E42").Value = "  +7.67%  "
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("E50").Value = "  -3.70%  "
$ws.Range("E51").Value = "  +0.01%  "
